$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (F1) to the new headers
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Updated metric values for row 2
$ws.Range("B2").Value = 0.05152593465593191
$ws.Range("C2").Value = 0.9984847101675268
$ws.Range("D2").Value = 0.1669273891272094

# New data cells
$ws.Range("G2").Value = 0.1260932844166139
$ws.Range("H2").Value = 0.991
